# Refresh the cryptos table (prices / 1h volume %) to the latest scrape.
# Source values are plain text (e.g. "562.35", "69.241.46" -- note the
# thousands-grouped prices use "." as a separator, so they are NOT valid
# numbers) -- a leading apostrophe forces Excel to keep them as text,
# matching the original inlineStr cells instead of auto-converting to Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '69.241.46'
$ws.Range("E2").Value = '  +2.85%  '

$ws.Range("D3").Value = "'" + '2.418.19'
$ws.Range("E3").Value = '  +1.72%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = "'" + '562.35'
$ws.Range("E5").Value = '  +2.84%  '

$ws.Range("D6").Value = "'" + '165.92'
$ws.Range("E6").Value = '  +6.39%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("D8").Value = "'" + '0.512'
$ws.Range("E8").Value = '  +2.40%  '

$ws.Range("D9").Value = "'" + '0.168'
$ws.Range("E9").Value = '  +8.88%  '

$ws.Range("D10").Value = "'" + '2.416.02'
$ws.Range("E10").Value = '  +1.38%  '

$ws.Range("E11").Value = '  -1.71%  '

$ws.Range("E12").Value = '  +3.04%  '

$ws.Range("E13").Value = '  -1.09%  '

$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = "'" + '69.241.19'
$ws.Range("E14").Value = '  +3.01%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = "'" + '0.0000177'
$ws.Range("E15").Value = '  +6.47%  '

$ws.Range("D16").Value = "'" + '2.862.93'
$ws.Range("E16").Value = '  -1.47%  '

$ws.Range("D17").Value = "'" + '23.91'
$ws.Range("E17").Value = '  +5.93%  '

$ws.Range("D18").Value = "'" + '2.422.08'
$ws.Range("E18").Value = '  +2.30%  '

$ws.Range("D19").Value = "'" + '10.79'
$ws.Range("E19").Value = '  +5.66%  '

$ws.Range("D20").Value = "'" + '342.11'
$ws.Range("E20").Value = '  +4.88%  '

$ws.Range("D21").Value = "'" + '7.11'
$ws.Range("E21").Value = '  +6.21%  '

$ws.Range("E22").Value = '  +3.88%  '

$ws.Range("E23").Value = '  +8.36%  '

$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("D25").Value = "'" + '65.95'
$ws.Range("E25").Value = '  +1.28%  '

$ws.Range("D26").Value = "'" + '3.80'
$ws.Range("E26").Value = '  +6.97%  '

$ws.Range("D27").Value = "'" + '8.46'
$ws.Range("E27").Value = '  +7.48%  '

$ws.Range("D28").Value = "'" + '2.543.30'
$ws.Range("E28").Value = '  +1.92%  '

$ws.Range("E29").Value = '  +0.49%  '

$ws.Range("E30").Value = '  +8.31%  '

$ws.Range("D31").Value = "'" + '7.37'
$ws.Range("E31").Value = '  +6.68%  '

$ws.Range("E32").Value = '  +12.26%  '

$ws.Range("D33").Value = "'" + '453.34'
$ws.Range("E33").Value = '  +10.66%  '

$ws.Range("E34").Value = '  +0.07%  '

$ws.Range("E35").Value = '  +2.39%  '

$ws.Range("D36").Value = "'" + '158.30'
$ws.Range("E36").Value = '  -0.84%  '

$ws.Range("D37").Value = "'" + '19.08'
$ws.Range("E37").Value = '  +0.77%  '

$ws.Range("E38").Value = '  +0.05%  '

$ws.Range("E39").Value = '  +5.98%  '

$ws.Range("D40").Value = "'" + '18.19'
$ws.Range("E40").Value = '  +3.76%  '

$ws.Range("E41").Value = '  +4.50%  '

$ws.Range("E42").Value = '  +5.48%  '

$ws.Range("E43").Value = '  +5.45%  '

$ws.Range("D44").Value = "'" + '37.77'
$ws.Range("E44").Value = '  +1.85%  '

$ws.Range("E45").Value = '  +4.27%  '

$ws.Range("E46").Value = '  +9.46%  '

$ws.Range("D47").Value = "'" + '134.68'
$ws.Range("E47").Value = '  +5.90%  '

$ws.Range("E48").Value = '  +4.70%  '

$ws.Range("D49").Value = "'" + '0.0722'
$ws.Range("E49").Value = '  +3.13%  '

$ws.Range("E50").Value = '  +4.73%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = "'" + '0.559'
$ws.Range("E51").Value = '  +2.37%  '
